# Domino Single Ethernet Rev. C: changed Logos
# Rename the worksheet from Rev. B to Rev. C (this also updates the
# sheet-qualified references inside the Print_Area defined names).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Domino Single Ethernet Rev. C"

# Re-assert the plain _xlnm.Print_Area defined name so that it also
# points at the renamed sheet with the original $A$1:$I$14 range.
$ws.PageSetup.PrintArea = '$A$1:$I$14'

# Slightly widened column widths (as recalculated when the sheet was
# resaved after the logo artwork changed).
$ws.Columns.Item(1).ColumnWidth = 4.166666666666667
$ws.Columns.Item(2).ColumnWidth = 4.166666666666667
$ws.Columns.Item(3).ColumnWidth = 26.166666666666668
$ws.Columns.Item(4).ColumnWidth = 28.166666666666668
$ws.Columns.Item(5).ColumnWidth = 31.333333333333332
$ws.Columns.Item(6).ColumnWidth = 26.166666666666668
$ws.Columns.Item(7).ColumnWidth = 40.5
$ws.Columns.Item(8).ColumnWidth = 63.0
$ws.Columns.Item(9).ColumnWidth = 24.5
